$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "73.266.32"
$ws.Range("E2").Value = "  +2.09%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "4.049.02"
$ws.Range("E3").Value = "  +1.26%  "

$ws.Range("E4").Value = "  -0.06%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "581.09"
$ws.Range("E5").Value = "  +9.61%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "152.01"
$ws.Range("E6").Value = "  +0.94%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "4.044.39"
$ws.Range("E7").Value = "  +1.34%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.691"
$ws.Range("E8").Value = "  +0.07%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.758"
$ws.Range("E10").Value = "  +2.18%  "

$ws.Range("E11").Value = "  -0.23%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "53.41"
$ws.Range("E12").Value = "  +12.69%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000324"
$ws.Range("E13").Value = "  -1.15%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "11.10"
$ws.Range("E14").Value = "  +4.69%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.702.18"
$ws.Range("E15").Value = "  +1.20%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "4.057.01"
$ws.Range("E16").Value = "  +1.31%  "

$ws.Range("E17").Value = "  +2.69%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("B18").Value = "Polygon"
$ws.Range("C18").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D18").Value = "1.23"
$ws.Range("E18").Value = "  +3.58%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("B19").Value = "Chainlink"
$ws.Range("C19").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D19").Value = "20.73"
$ws.Range("E19").Value = "  +1.27%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "73.245.54"
$ws.Range("E20").Value = "  +2.16%  "

$ws.Range("E21").Value = "  -0.36%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "442.48"
$ws.Range("E22").Value = "  +4.07%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("B23").Value = "PancakeSwap"
$ws.Range("C23").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D23").Value = "4.62"
$ws.Range("E23").Value = "  +10.39%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("B24").Value = "Litecoin"
$ws.Range("C24").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D24").Value = "97.78"
$ws.Range("E24").Value = "  +0.48%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.54"

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "14.56"
$ws.Range("E26").Value = "  +1.75%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "4.29"
$ws.Range("E27").Value = "  +19.48%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "11.57"
$ws.Range("E28").Value = "  +3.10%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "10.94"
$ws.Range("E29").Value = "  +2.72%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.98"
$ws.Range("E30").Value = "  +2.54%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "37.02"
$ws.Range("E31").Value = "  +1.38%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.92"
$ws.Range("E32").Value = "  +13.41%  "

$ws.Range("E33").Value = "  +4.02%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "13.66"
$ws.Range("E34").Value = "  +2.58%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "693.10"

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "48.48"
$ws.Range("E36").Value = "  +10.43%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "67.83"
$ws.Range("E37").Value = "  +3.37%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("B38").Value = "TheGraph"
$ws.Range("C38").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D38").Value = "0.447"
$ws.Range("E38").Value = "  +3.48%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("B39").Value = "PEPE"
$ws.Range("C39").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D39").Value = "0.0₃0890"
$ws.Range("E39").Value = "  +8.13%  "

$ws.Range("E40").Value = "  -1.33%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "11.31"
$ws.Range("E41").Value = "  +18.14%  "

$ws.Range("E42").Value = "  -1.14%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.999"
$ws.Range("E43").Value = "  +0.05%  "

$ws.Range("E44").Value = "  +4.84%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0496"
$ws.Range("E45").Value = "  +2.25%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.00"
$ws.Range("E46").Value = "  +0.07%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("B47").Value = "Fetch.AI"
$ws.Range("C47").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D47").Value = "2.76"
$ws.Range("E47").Value = "  +5.48%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("B48").Value = "Stellar"
$ws.Range("C48").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D48").Value = "0.151"
$ws.Range("E48").Value = "  +0.91%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("B49").Value = "ApeXProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D49").Value = "3.40"
$ws.Range("E49").Value = "  -0.87%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("B50").Value = "LidoDAOToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D50").Value = "3.49"
$ws.Range("E50").Value = "  +6.57%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("B51").Value = "Stacks"
$ws.Range("C51").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D51").Value = "3.06"
$ws.Range("E51").Value = "  +2.91%  "

